$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Copy style from K1 (existing header) to L1:N1
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats

# Data values
$ws.Range("L2").Value = 92.19084927254757
$ws.Range("M2").Value = 238184
$ws.Range("N2").Value = 321.4358974358975

$ws.Range("L3").Value = 82.97112949625762
$ws.Range("M3").Value = 25369
$ws.Range("N3").Value = 285.0449438202247

$ws.Range("L4").Value = 90.90970086312072
$ws.Range("M4").Value = 186820
$ws.Range("N4").Value = 145.953125

$ws.Range("L5").Value = 81.91974458046401
$ws.Range("M5").Value = 16826
$ws.Range("N5").Value = 163.3592233009709

$ws.Range("L6").Value = 19.81586961574516
$ws.Range("M6").Value = 2101
$ws.Range("N6").Value = 15.56296296296296

$ws.Range("L7").Value = 14.52164198228234
$ws.Range("M7").Value = 107
$ws.Range("N7").Value = 6.294117647058823
